$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Inscritos" (E) and related totals for a handful of rows
# Row 21
$ws.Range("E21").Value = 122

# Row 26
$ws.Range("E26").Value = 113

# Row 30
$ws.Range("E30").Value = 167
$ws.Range("F30").Value = 93
$ws.Range("H30").Value = 93

# Row 33
$ws.Range("E33").Value = 236

# Row 36
$ws.Range("E36").Value = 50

# Row 38
$ws.Range("E38").Value = 79

# Row 39
$ws.Range("E39").Value = 157

# Row 42
$ws.Range("E42").Value = 281

# Row 48
$ws.Range("E48").Value = 164
